$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushes the old row 8..20 down to 9..21),
# matching the diff's new dimension A1:R21.
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C8").Value = 'Ñuble'
$ws.Range("D8").Value = 44498
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112022
$ws.Range("G8").Value = 'Arveja Verde'
$ws.Range("H8").Value = 'Sin especificar'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17500
$ws.Range("N8").Value = '$/saco 25 kilos'
$ws.Range("O8").Value = 'Región del Maule'
$ws.Range("P8").Value = 700
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = 'Hortaliza'
